$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The edit reorders several paragraphs' worth of content around the
# "Objetivos" / "Docente(s) Responsavel(eis)" / "Programa resumido" /
# "Programa" / "Avaliacao" / "Bibliografia" sections. Paragraph *styles*
# and run *formatting* (bold/italic) all stay exactly where they are;
# only the plain text content moves between paragraphs (and, in one case,
# between runs inside the same paragraph). So the safest approach is to
# capture the original text of every paragraph/run we need (by Word
# paragraph index, 1-based) BEFORE mutating anything, then write the
# captured strings back into their new homes.
# ---------------------------------------------------------------------------

# --- 1. Capture source text blocks (before any mutation) -------------------

# Paragraph 6: "Conferir aos alunos ..." (Objetivos, PT)
$txt_objetivosPT = $d.Paragraphs.Item(6).Range.Text
$txt_objetivosPT = $txt_objetivosPT.Substring(0, $txt_objetivosPT.Length - 1)

# Paragraph 7: "Providing to the students ..." (Objetivos, EN, italic)
$txt_objetivosEN = $d.Paragraphs.Item(7).Range.Text
$txt_objetivosEN = $txt_objetivosEN.Substring(0, $txt_objetivosEN.Length - 1)

# Paragraph 11: "Processo Quimico e Industria Quimica ..." (Programa resumido, PT)
$txt_resumidoPT = $d.Paragraphs.Item(11).Range.Text
$txt_resumidoPT = $txt_resumidoPT.Substring(0, $txt_resumidoPT.Length - 1)

# Paragraph 12: "Chemical Process and Chemical Industry ..." (Programa resumido, EN, italic)
$txt_resumidoEN = $d.Paragraphs.Item(12).Range.Text
$txt_resumidoEN = $txt_resumidoEN.Substring(0, $txt_resumidoEN.Length - 1)

# Paragraph 14: "O conteudo desta disciplina sera ..." (Programa, PT)
$txt_programaPT = $d.Paragraphs.Item(14).Range.Text
$txt_programaPT = $txt_programaPT.Substring(0, $txt_programaPT.Length - 1)

# Paragraph 9: "1285870 - Marcos Villela Barcza" (Docente responsavel)
$txt_docente = $d.Paragraphs.Item(9).Range.Text
$txt_docente = $txt_docente.Substring(0, $txt_docente.Length - 1)

# Paragraph 19: bibliography block (multiple runs, many <w:br/>)
$txt_bibliografia = $d.Paragraphs.Item(19).Range.Text
$txt_bibliografia = $txt_bibliografia.Substring(0, $txt_bibliografia.Length - 1)

# Paragraph 17 (Avaliacao) internal content blocks, located relative to
# their bold labels so we don't depend on exact character offsets.
$p17 = $d.Paragraphs.Item(17)

$metodoLabel = $p17.Range.Duplicate
$metodoLabel.Find.Execute("Método: ")

$criterioLabel = $p17.Range.Duplicate
$criterioLabel.Find.Execute("Critério: ")

$normaLabel = $p17.Range.Duplicate
$normaLabel.Find.Execute("Norma de recuperação: ")

$run_metodo = $d.Range($metodoLabel.End, $criterioLabel.Start)
$txt_metodoBody = $run_metodo.Text
# Drop the trailing line-break mark (<w:br/> -> chr(11)): in its original
# home it separated "Método:" from "Critério:" inside one paragraph; once
# it becomes its own standalone paragraph (14) that break is not wanted.
if ($txt_metodoBody.Length -gt 0 -and [int][char]$txt_metodoBody.Substring($txt_metodoBody.Length - 1, 1) -eq 11) {
    $txt_metodoBody = $txt_metodoBody.Substring(0, $txt_metodoBody.Length - 1)
}

$run_criterio = $d.Range($criterioLabel.End, $normaLabel.Start)
$txt_criterioBody = $run_criterio.Text

$run_norma = $d.Range($normaLabel.End, $p17.Range.End - 1)
$txt_normaBody = $run_norma.Text

# --- 2. Write the captured text back into its new home ---------------------

# Paragraph 6 becomes the "Programa resumido" PT text.
$d.Paragraphs.Item(6).Range.Find.Execute($txt_objetivosPT, $true, $false, $false, $false, $false, $true, 1, $false, $txt_resumidoPT, 2)

# Paragraph 7 becomes the "Programa resumido" EN text (keeps italic run).
$d.Paragraphs.Item(7).Range.Find.Execute($txt_objetivosEN, $true, $false, $false, $false, $false, $true, 1, $false, $txt_resumidoEN, 2)

# Paragraph 9 (List Bullet) becomes the "Objetivos" PT text.
$d.Paragraphs.Item(9).Range.Find.Execute($txt_docente, $true, $false, $false, $false, $false, $true, 1, $false, $txt_objetivosPT, 2)

# Paragraph 11 becomes the "Programa" PT text.
$d.Paragraphs.Item(11).Range.Find.Execute($txt_resumidoPT, $true, $false, $false, $false, $false, $true, 1, $false, $txt_programaPT, 2)

# Paragraph 12 becomes the "Objetivos" EN text (keeps italic run).
$d.Paragraphs.Item(12).Range.Find.Execute($txt_resumidoEN, $true, $false, $false, $false, $false, $true, 1, $false, $txt_objetivosEN, 2)

# Paragraph 14 becomes the "Metodo" body text (was inside paragraph 17).
$d.Paragraphs.Item(14).Range.Find.Execute($txt_programaPT, $true, $false, $false, $false, $false, $true, 1, $false, $txt_metodoBody, 2)

# Paragraph 19 becomes the "Docente responsavel" text.
$d.Paragraphs.Item(19).Range.Find.Execute($txt_bibliografia, $true, $false, $false, $false, $false, $true, 1, $false, $txt_docente, 2)

# --- 3. Rearrange the three content runs inside paragraph 17 ---------------
# Re-resolve the label ranges (the document hasn't shifted positions for
# paragraph 17 itself yet, but re-find to be safe) and rewrite each body
# in one pass, from the end of the paragraph backwards so earlier offsets
# stay valid.

$p17 = $d.Paragraphs.Item(17)

$metodoLabel = $p17.Range.Duplicate
$metodoLabel.Find.Execute("Método: ")

$criterioLabel = $p17.Range.Duplicate
$criterioLabel.Find.Execute("Critério: ")

$normaLabel = $p17.Range.Duplicate
$normaLabel.Find.Execute("Norma de recuperação: ")

$run_norma = $d.Range($normaLabel.End, $p17.Range.End - 1)
$run_norma.Text = $txt_bibliografia

$run_criterio = $d.Range($criterioLabel.End, $normaLabel.Start)
$run_criterio.Text = $txt_normaBody + [char]11

$run_metodo = $d.Range($metodoLabel.End, $criterioLabel.Start)
$run_metodo.Text = $txt_criterioBody

Write-Output "done"
